# Scheduled market-data refresh: updates computed leve-profit columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Table_<job>" sheets with
# freshly pulled Universalis market prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2229.1428
$ws.Range("I40").Value = 2490.1
$ws.Range("K40").Value = 2490.1
$ws.Range("M40").Value = -2315.1
$ws.Range("H98").Value = 1190.44
$ws.Range("I98").Value = 1163.5217
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1163.5217
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 334.4783
$ws.Range("N98").Value = -4496
$ws.Range("H113").Value = 3432.5789
$ws.Range("I113").Value = 3233.875
$ws.Range("J113").Value = 3577.0908
$ws.Range("K113").Value = 3233.875
$ws.Range("L113").Value = 3577.0908
$ws.Range("M113").Value = 20.125
$ws.Range("N113").Value = -10085.0908
$ws.Range("H122").Value = 1190.44
$ws.Range("I122").Value = 1163.5217
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3490.5651
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1040.5651
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 2149.932
$ws.Range("I132").Value = 1737.4849
$ws.Range("J132").Value = 3387.2727
$ws.Range("K132").Value = 5212.4547
$ws.Range("L132").Value = 10161.8181
$ws.Range("M132").Value = -2682.4547
$ws.Range("N132").Value = -15221.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7009.5947
$ws.Range("I32").Value = 5596.73
$ws.Range("J32").Value = 12572.75
$ws.Range("K32").Value = 5596.73
$ws.Range("L32").Value = 12572.75
$ws.Range("M32").Value = -5309.73
$ws.Range("N32").Value = -13146.75
$ws.Range("H45").Value = 1517188.2
$ws.Range("I45").Value = 2165915.8
$ws.Range("J45").Value = 3490.6667
$ws.Range("K45").Value = 2165915.8
$ws.Range("L45").Value = 3490.6667
$ws.Range("M45").Value = -2165538.8
$ws.Range("N45").Value = -4244.6667
$ws.Range("H122").Value = 5002863
$ws.Range("J122").Value = 12502000
$ws.Range("L122").Value = 37506000
$ws.Range("N122").Value = -37510900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7753906.5
$ws.Range("I86").Value = 8773949
$ws.Range("J86").Value = 1581.4
$ws.Range("K86").Value = 8773949
$ws.Range("L86").Value = 1581.4
$ws.Range("M86").Value = -8772826
$ws.Range("N86").Value = -3827.4
$ws.Range("H89").Value = 7753906.5
$ws.Range("I89").Value = 8773949
$ws.Range("J89").Value = 1581.4
$ws.Range("K89").Value = 43869745
$ws.Range("L89").Value = 7907
$ws.Range("M89").Value = -43864129
$ws.Range("N89").Value = -19139
$ws.Range("H94").Value = 1016.7778
$ws.Range("I94").Value = 845.7826
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 845.7826
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -394.7826
$ws.Range("N94").Value = -2902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1607.7222
$ws.Range("I99").Value = 1710.1818
$ws.Range("J99").Value = 1446.7142
$ws.Range("K99").Value = 1710.1818
$ws.Range("L99").Value = 1446.7142
$ws.Range("M99").Value = -212.1818000000001
$ws.Range("N99").Value = -4442.7142
$ws.Range("H122").Value = 10864.259
$ws.Range("I122").Value = 3471.8948
$ws.Range("J122").Value = 28421.125
$ws.Range("K122").Value = 10415.6844
$ws.Range("L122").Value = 85263.375
$ws.Range("M122").Value = -7965.6844
$ws.Range("N122").Value = -90163.375
$ws.Range("H126").Value = 1607.7222
$ws.Range("I126").Value = 1710.1818
$ws.Range("J126").Value = 1446.7142
$ws.Range("K126").Value = 5130.5454
$ws.Range("L126").Value = 4340.142599999999
$ws.Range("M126").Value = -2660.5454
$ws.Range("N126").Value = -9280.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H57").Value = 2440
$ws.Range("I57").Value = 800
$ws.Range("J57").Value = 2850
$ws.Range("K57").Value = 2400
$ws.Range("L57").Value = 8550
$ws.Range("M57").Value = -1841
$ws.Range("N57").Value = -9668
$ws.Range("H64").Value = 2436.3103
$ws.Range("I64").Value = 1125.6666
$ws.Range("J64").Value = 2778.2173
$ws.Range("K64").Value = 3376.9998
$ws.Range("L64").Value = 8334.651899999999
$ws.Range("M64").Value = -3106.9998
$ws.Range("N64").Value = -8874.651899999999
$ws.Range("H67").Value = 2436.3103
$ws.Range("I67").Value = 1125.6666
$ws.Range("J67").Value = 2778.2173
$ws.Range("K67").Value = 3376.9998
$ws.Range("L67").Value = 8334.651899999999
$ws.Range("M67").Value = -2440.9998
$ws.Range("N67").Value = -10206.6519
$ws.Range("H97").Value = 8831.429
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 14705
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 44115
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -45107
$ws.Range("H114").Value = 428.14285
$ws.Range("I114").Value = 291.27274
$ws.Range("J114").Value = 930
$ws.Range("K114").Value = 873.81822
$ws.Range("L114").Value = 2790
$ws.Range("M114").Value = 2380.18178
$ws.Range("N114").Value = -9298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H70").Value = 5446.231
$ws.Range("I70").Value = 5048.75
$ws.Range("J70").Value = 5786.9287
$ws.Range("K70").Value = 5048.75
$ws.Range("L70").Value = 5786.9287
$ws.Range("M70").Value = -4778.75
$ws.Range("N70").Value = -6326.9287
$ws.Range("H73").Value = 5446.231
$ws.Range("I73").Value = 5048.75
$ws.Range("J73").Value = 5786.9287
$ws.Range("K73").Value = 5048.75
$ws.Range("L73").Value = 5786.9287
$ws.Range("M73").Value = -4112.75
$ws.Range("N73").Value = -7658.9287
$ws.Range("H102").Value = 3651.7273
$ws.Range("I102").Value = 4037.4285
$ws.Range("J102").Value = 2976.75
$ws.Range("K102").Value = 4037.4285
$ws.Range("L102").Value = 2976.75
$ws.Range("M102").Value = -2415.4285
$ws.Range("N102").Value = -6220.75
$ws.Range("H107").Value = 282.30768
$ws.Range("I107").Value = 105.875
$ws.Range("J107").Value = 564.6
$ws.Range("K107").Value = 105.875
$ws.Range("L107").Value = 564.6
$ws.Range("M107").Value = 1814.125
$ws.Range("N107").Value = -4404.6
$ws.Range("H126").Value = 3124.2856
$ws.Range("I126").Value = 2265
$ws.Range("J126").Value = 4842.857
$ws.Range("K126").Value = 6795
$ws.Range("L126").Value = 14528.571
$ws.Range("M126").Value = -4325
$ws.Range("N126").Value = -19468.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5204.3184
$ws.Range("I7").Value = 3325.4167
$ws.Range("J7").Value = 7459
$ws.Range("K7").Value = 3325.4167
$ws.Range("L7").Value = 7459
$ws.Range("M7").Value = -3213.4167
$ws.Range("N7").Value = -7683
$ws.Range("H40").Value = 3831.2415
$ws.Range("I40").Value = 3550.0454
$ws.Range("J40").Value = 4715
$ws.Range("K40").Value = 3550.0454
$ws.Range("L40").Value = 4715
$ws.Range("M40").Value = -3414.0454
$ws.Range("N40").Value = -4987
$ws.Range("H46").Value = 756.25
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 858.3333
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 858.3333
$ws.Range("M46").Value = -262
$ws.Range("N46").Value = -1234.3333
$ws.Range("H100").Value = 3426
$ws.Range("I100").Value = 1998.6471
$ws.Range("J100").Value = 6122.1113
$ws.Range("K100").Value = 1998.6471
$ws.Range("L100").Value = 6122.1113
$ws.Range("M100").Value = -1457.6471
$ws.Range("N100").Value = -7204.1113
$ws.Range("H122").Value = 4327.3613
$ws.Range("I122").Value = 3831.6
$ws.Range("J122").Value = 5454.091
$ws.Range("K122").Value = 11494.8
$ws.Range("L122").Value = 16362.273
$ws.Range("M122").Value = -9044.799999999999
$ws.Range("N122").Value = -21262.273
$ws.Range("H126").Value = 5204.3184
$ws.Range("I126").Value = 3325.4167
$ws.Range("J126").Value = 7459
$ws.Range("K126").Value = 9976.250100000001
$ws.Range("L126").Value = 22377
$ws.Range("M126").Value = -7506.250100000001
$ws.Range("N126").Value = -27317
$ws.Range("H132").Value = 3726.8774
$ws.Range("I132").Value = 3490.6
$ws.Range("J132").Value = 4317.5713
$ws.Range("K132").Value = 10471.8
$ws.Range("L132").Value = 12952.7139
$ws.Range("M132").Value = -7941.799999999999
$ws.Range("N132").Value = -18012.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 13024
$ws.Range("I17").Value = 5060
$ws.Range("J17").Value = 18333.334
$ws.Range("K17").Value = 5060
$ws.Range("L17").Value = 18333.334
$ws.Range("M17").Value = -4888
$ws.Range("N17").Value = -18677.334
$ws.Range("H41").Value = 9550.6
$ws.Range("J41").Value = 9550.6
$ws.Range("L41").Value = 9550.6
$ws.Range("N41").Value = -10330.6
$ws.Range("H45").Value = 5480
$ws.Range("I45").Value = 6000
$ws.Range("J45").Value = 5350
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 5350
$ws.Range("M45").Value = -5509
$ws.Range("N45").Value = -6332
$ws.Range("H74").Value = 11037.333
$ws.Range("J74").Value = 11037.333
$ws.Range("L74").Value = 11037.333
$ws.Range("N74").Value = -12909.333
$ws.Range("H77").Value = 11037.333
$ws.Range("J77").Value = 11037.333
$ws.Range("L77").Value = 33111.999
$ws.Range("N77").Value = -42471.999

Write-Output "Updated 248 cells across $($wb.Worksheets.Count) sheets"
